$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the two new boundary-condition data rows (0 < x < 200 triangle classification)
$ws.Range("A5").Value = "0,4,5"
$ws.Range("B5").Value = "Not a Triangle"
$ws.Range("A6").Value = "200,4,5"
$ws.Range("B6").Value = "Not a Triangle"

# Build a single combined center/middle alignment style and seed it on A1,
# then fan it out to the whole table via a format-only paste so every
# cell ends up sharing the very same style record.
$centerStyle = $wb.Styles.Add("CenterMiddle")
$centerStyle.HorizontalAlignment = -4108
$centerStyle.VerticalAlignment = -4108
Write-Output "style ready"
$ws.Range("A1").Style = "CenterMiddle"
$ws.Range("A1").Copy()
$ws.Range("A1:B6").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Column widths to fit the (now wider) content
$ws.Columns.Item(1).ColumnWidth = 11.428571428571429
$ws.Columns.Item(2).ColumnWidth = 18.714285714285715

# Move the active selection to A2
$ws.Range("A2").Select()
